# The commit swaps the colour scheme of the presentation's theme
# (ppt/theme/theme1.xml, the theme used by the slide master / all the
# slides) from the "Integral" palette to the stock Office 2013+ default
# palette ("Office Theme"). The font scheme (Arial-based "Office" font
# scheme) and the format scheme (fills/lines/effects) are identical
# between the two themes already, so only the 12 colour-scheme slots
# need to change.
#
# PowerPoint's ThemeColorScheme.Item(n).RGB uses the classic COM/OLE
# "COLORREF" packing: R + G*256 + B*65536 (i.e. 0x00BBGGRR), so convert
# each target hex colour with a small helper instead of hand-swapping
# bytes everywhere.

function ToComRgb($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

function HexToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ToComRgb $r $g $b
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Target palette: the default Office theme colours.
$tcs.Item(1).RGB  = HexToComRgb "000000"   # dk1
$tcs.Item(2).RGB  = HexToComRgb "FFFFFF"   # lt1
$tcs.Item(3).RGB  = HexToComRgb "44546A"   # dk2
$tcs.Item(4).RGB  = HexToComRgb "E7E6E6"   # lt2
$tcs.Item(5).RGB  = HexToComRgb "5B9BD5"   # accent1
$tcs.Item(6).RGB  = HexToComRgb "ED7D31"   # accent2
$tcs.Item(7).RGB  = HexToComRgb "A5A5A5"   # accent3
$tcs.Item(8).RGB  = HexToComRgb "FFC000"   # accent4
$tcs.Item(9).RGB  = HexToComRgb "4472C4"   # accent5
$tcs.Item(10).RGB = HexToComRgb "70AD47"   # accent6
$tcs.Item(11).RGB = HexToComRgb "0563C1"   # hlink
$tcs.Item(12).RGB = HexToComRgb "954F72"   # folHlink
